# Fix dataset problems in "num_commenti_hate_per_topic" sheet.
# - Row 3's topic had a typo introduced: CRONACA -> CROANCA
# - All topic/social rows shifted down by one row (a new POLITICA/YouTube
#   row was appended at the bottom, row 12)
# - All the num_hate_topic (column D) counts were corrected/updated

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New topic/social/count data for rows 3-12 (columns B, C, D)
$data = @(
    @("CROANCA",      "Facebook",  7),
    @("CRONACA",       "Facebook",  94),
    @("CRONACA",       "Instagram", 198),
    @("CRONACA",       "YouTube",   88),
    @("CRONACA NERA",  "Facebook",  186),
    @("CRONACA NERA",  "Instagram", 204),
    @("CRONACA NERA",  "YouTube",   83),
    @("POLITICA",      "Facebook",  88),
    @("POLITICA",      "Instagram", 192),
    @("POLITICA",      "YouTube",   91)
)

$row = 3
foreach ($entry in $data) {
    $ws.Cells.Item($row, 2).Value = $entry[0]
    $ws.Cells.Item($row, 3).Value = $entry[1]
    $ws.Cells.Item($row, 4).Value = $entry[2]
    $row = $row + 1
}

# The newly appended row (12) needs an (empty) column-A cell too, matching
# the pattern already used by rows 3-11 in this table. Re-applying the
# "Normal" style forces Excel to materialize the otherwise-empty cell
# without altering any formatting.
$ws.Cells.Item(12, 1).Style = "Normal"
